# Apply the SAO_JOSE_DO_NORTE.xlsx update:
#  - Remove the "Desarquivamentos Pendentes" sheet (no longer used)
#  - Rename "Paineis DARQ" -> "PAINEIS DARQ"
#  - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"

$wb = $excel.ActiveWorkbook

# Avoid the "are you sure you want to delete" confirmation dialog.
$excel.DisplayAlerts = $false

$wb.Worksheets("Desarquivamentos Pendentes").Delete()

$wb.Worksheets("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $true
